$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns B-E stay text-typed so values like percentages and
# numeric-looking strings are preserved exactly as literal text,
# matching the original inline-string cell contents.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.84%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.35"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.66%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "9.00%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07997"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.34%"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.984"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.34%"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.612"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.23%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.947"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.36%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9262"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.12%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1248"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-5.89%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1948"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.94%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.714"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "24.72%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09136"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.14%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03634"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.76%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1050"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "9.59%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001311"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.99%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006189"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "3.23%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.350"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.44%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.510"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.35%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.63%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.43%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-4.37%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04410"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.01%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001264"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.37%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004588"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.07%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001151"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-3.34%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02527"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.68%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05381"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.85%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007426"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.63%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009514"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.23%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.64%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002119"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.02%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01073"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.88%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006783"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.91%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.04%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002973"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-11.10%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002293"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-7.64%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.04%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.04%"
